# Change the table style on the B1/B2 table (Slide 5, Shape 2) from the
# custom "Table_0" style ({ABE6603E-14B3-4F9A-B1EB-FC1B370FFF31}) to the
# built-in "No Style, Table Grid" style ({85A7E0DF-00E8-4060-950A-1A658CDCE999}).

$p = $ppt.ActivePresentation

$targetSlideIndex = 0
$targetShapeIndex = 0

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shp = $s.Shapes.Item($j)
        if ($shp.HasTable) {
            $targetSlideIndex = $i
            $targetShapeIndex = $j
        }
    }
}

$slide = $p.Slides.Item($targetSlideIndex)
$shape = $slide.Shapes.Item($targetShapeIndex)
$table = $shape.Table

$table.ApplyStyle("{85A7E0DF-00E8-4060-950A-1A658CDCE999}")
